$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2830809886899601
$ws.Cells.Item(2, 3).Value = 0.0522837017491895
$ws.Cells.Item(2, 4).Value = 0.0330633798573885
$ws.Cells.Item(2, 5).Value = 0.1652192017219249
$ws.Cells.Item(2, 6).Value = 0.8129378803801544
$ws.Cells.Item(2, 11).Value = 0.2587585588637467
$ws.Cells.Item(2, 13).Value = 0.2176760785492249
$ws.Cells.Item(2, 15).Value = 2.832204797704463
$ws.Cells.Item(3, 2).Value = 0.2507523373220124
$ws.Cells.Item(3, 3).Value = 0.04908269788906239
$ws.Cells.Item(3, 4).Value = 0.03135415875789249
$ws.Cells.Item(3, 5).Value = 0.1541697129981685
$ws.Cells.Item(3, 6).Value = 0.8105211648305968
$ws.Cells.Item(3, 11).Value = 0.2259095237993876
$ws.Cells.Item(3, 13).Value = 0.1957731359861299
$ws.Cells.Item(3, 15).Value = 2.838189110713785
$ws.Cells.Item(4, 2).Value = 0.2309166325486842
$ws.Cells.Item(4, 3).Value = 0.04710047589718158
$ws.Cells.Item(4, 4).Value = 0.03029392869182601
$ws.Cells.Item(4, 5).Value = 0.1474994569335948
$ws.Cells.Item(4, 6).Value = 0.8095404201753809
$ws.Cells.Item(4, 11).Value = 0.205696411878975
$ws.Cells.Item(4, 13).Value = 0.1823977817767144
$ws.Cells.Item(4, 15).Value = 2.843528099202757
$ws.Cells.Item(5, 2).Value = 0.2228373722527976
$ws.Cells.Item(5, 3).Value = 0.04628852206153056
$ws.Cells.Item(5, 4).Value = 0.02985919589619357
$ws.Cells.Item(5, 5).Value = 0.1448097466896172
$ws.Cells.Item(5, 6).Value = 0.8092671998140943
$ws.Cells.Item(5, 11).Value = 0.1974487944272738
$ws.Cells.Item(5, 13).Value = 0.176965547393074
$ws.Cells.Item(5, 15).Value = 2.846122018784854
$ws.Cells.Item(6, 2).Value = 0.2214960669335255
$ws.Cells.Item(6, 3).Value = 0.0461534463431903
$ws.Cells.Item(6, 4).Value = 0.02978684771185414
$ws.Cells.Item(6, 5).Value = 0.1443648351530982
$ws.Cells.Item(6, 6).Value = 0.8092294663289721
$ws.Cells.Item(6, 11).Value = 0.1960786539932968
$ws.Cells.Item(6, 13).Value = 0.1760646344298564
$ws.Cells.Item(6, 15).Value = 2.846577987859177
$ws.Cells.Item(7, 2).Value = 0.2308076563246289
$ws.Cells.Item(7, 3).Value = 0.04708954248439312
$ws.Cells.Item(7, 4).Value = 0.03028807655272203
$ws.Cells.Item(7, 5).Value = 0.1474630675859743
$ws.Cells.Item(7, 6).Value = 0.8095362235761243
$ws.Cells.Item(7, 11).Value = 0.2055852240915073
$ws.Cells.Item(7, 13).Value = 0.1823244466209744
$ws.Cells.Item(7, 15).Value = 2.843561388733775
$ws.Cells.Item(8, 2).Value = 0.2719313333304285
$ws.Cells.Item(8, 3).Value = 0.05118350413928852
$ws.Cells.Item(8, 4).Value = 0.03247628878721542
$ws.Cells.Item(8, 5).Value = 0.1613854949032074
$ws.Cells.Item(8, 6).Value = 0.8120001246013402
$ws.Cells.Item(8, 11).Value = 0.2474414878396942
$ws.Cells.Item(8, 13).Value = 0.2101087198407185
$ws.Cells.Item(8, 15).Value = 2.833922492611748
$ws.Cells.Item(9, 2).Value = 0.3526751361916354
$ws.Cells.Item(9, 3).Value = 0.05907710337417882
$ws.Cells.Item(9, 4).Value = 0.03668104365467428
$ws.Cells.Item(9, 5).Value = 0.1896054350021572
$ws.Cells.Item(9, 6).Value = 0.8208287424903773
$ws.Cells.Item(9, 11).Value = 0.3291626809500769
$ws.Cells.Item(9, 13).Value = 0.2651795853471342
$ws.Cells.Item(9, 15).Value = 2.828246759254682
$ws.Cells.Item(10, 2).Value = 0.4120481324848413
$ws.Cells.Item(10, 3).Value = 0.06479310223576817
$ws.Cells.Item(10, 4).Value = 0.03971667308556448
$ws.Cells.Item(10, 5).Value = 0.2109173303031753
$ws.Cells.Item(10, 6).Value = 0.8297605130629222
$ws.Cells.Item(10, 11).Value = 0.3889738941032306
$ws.Cells.Item(10, 13).Value = 0.3060093878631278
$ws.Cells.Item(10, 15).Value = 2.832169520643475
$ws.Cells.Item(11, 2).Value = 0.4390675447407659
$ws.Cells.Item(11, 3).Value = 0.06737509638109884
$ws.Cells.Item(11, 4).Value = 0.04108583131483812
$ws.Cells.Item(11, 5).Value = 0.2207425410970814
$ws.Cells.Item(11, 6).Value = 0.8343568068720373
$ws.Cells.Item(11, 11).Value = 0.4161319071561991
$ws.Cells.Item(11, 13).Value = 0.3246669330914855
$ws.Cells.Item(11, 15).Value = 2.835718031893919
$ws.Cells.Item(12, 2).Value = 0.4493002885963335
$ws.Cells.Item(12, 3).Value = 0.0683501748729185
$ws.Cells.Item(12, 4).Value = 0.04160258225358149
$ws.Cells.Item(12, 5).Value = 0.2244821261931236
$ws.Cells.Item(12, 6).Value = 0.8361740954580483
$ws.Cells.Item(12, 11).Value = 0.4264084038767066
$ws.Cells.Item(12, 13).Value = 0.33174425484723
$ws.Cells.Item(12, 15).Value = 2.837315889454118
$ws.Cells.Item(13, 2).Value = 0.4470964426118087
$ws.Cells.Item(13, 3).Value = 0.0681402935614841
$ws.Cells.Item(13, 4).Value = 0.04149136761173366
$ws.Cells.Item(13, 5).Value = 0.2236758910059962
$ws.Cells.Item(13, 6).Value = 0.8357792938092814
$ws.Cells.Item(13, 11).Value = 0.4241955232466808
$ws.Cells.Item(13, 13).Value = 0.3302194873132791
$ws.Cells.Item(13, 15).Value = 2.836960453352447
$ws.Cells.Item(14, 2).Value = 0.4399093794446571
$ws.Cells.Item(14, 3).Value = 0.06745537036738369
$ws.Cells.Item(14, 4).Value = 0.0411283793473558
$ws.Cells.Item(14, 5).Value = 0.2210498169678417
$ws.Cells.Item(14, 6).Value = 0.8345047771042857
$ws.Cells.Item(14, 11).Value = 0.4169775164231453
$ws.Cells.Item(14, 13).Value = 0.3252489455074254
$ws.Cells.Item(14, 15).Value = 2.835844393224107
$ws.Cells.Item(15, 2).Value = 0.4355072241562254
$ws.Cells.Item(15, 3).Value = 0.06703548657704061
$ws.Cells.Item(15, 4).Value = 0.04090581379364266
$ws.Cells.Item(15, 5).Value = 0.2194437514313918
$ws.Cells.Item(15, 6).Value = 0.8337340997447029
$ws.Cells.Item(15, 11).Value = 0.4125552701514437
$ws.Cells.Item(15, 13).Value = 0.3222059244043152
$ws.Cells.Item(15, 15).Value = 2.835193880512804
$ws.Cells.Item(16, 2).Value = 0.4102825234807028
$ws.Cells.Item(16, 3).Value = 0.06462399204062308
$ws.Cells.Item(16, 4).Value = 0.03962695625602208
$ws.Cells.Item(16, 5).Value = 0.2102778735239568
$ws.Cells.Item(16, 6).Value = 0.8294708721264072
$ws.Cells.Item(16, 11).Value = 0.3871980050866455
$ws.Cells.Item(16, 13).Value = 0.304791765628778
$ws.Cells.Item(16, 15).Value = 2.831973152188709
$ws.Cells.Item(17, 2).Value = 0.3948103437952
$ws.Cells.Item(17, 3).Value = 0.06313991489960813
$ws.Cells.Item(17, 4).Value = 0.03883938548419508
$ws.Cells.Item(17, 5).Value = 0.2046884613820552
$ws.Cells.Item(17, 6).Value = 0.8269921536286091
$ws.Cells.Item(17, 11).Value = 0.3716289448530574
$ws.Cells.Item(17, 13).Value = 0.2941302615771164
$ws.Cells.Item(17, 15).Value = 2.830449468906778
$ws.Cells.Item(18, 2).Value = 0.3859121526270144
$ws.Cells.Item(18, 3).Value = 0.06228459848874479
$ws.Cells.Item(18, 4).Value = 0.03838528989454915
$ws.Cells.Item(18, 5).Value = 0.2014858309370808
$ws.Cells.Item(18, 6).Value = 0.8256166404156602
$ws.Cells.Item(18, 11).Value = 0.3626693126245186
$ws.Cells.Item(18, 13).Value = 0.2880059367447956
$ws.Cells.Item(18, 15).Value = 2.829739098452364
$ws.Cells.Item(19, 2).Value = 0.3828995617441535
$ws.Cells.Item(19, 3).Value = 0.06199470991160183
$ws.Cells.Item(19, 4).Value = 0.03823135187459314
$ws.Cells.Item(19, 5).Value = 0.2004035718830863
$ws.Cells.Item(19, 6).Value = 0.8251595306882678
$ws.Cells.Item(19, 11).Value = 0.3596349367105063
$ws.Cells.Item(19, 13).Value = 0.2859336990883463
$ws.Cells.Item(19, 15).Value = 2.829527078247366
$ws.Cells.Item(20, 2).Value = 0.3964572834735236
$ws.Cells.Item(20, 3).Value = 0.06329807520505426
$ws.Cells.Item(20, 4).Value = 0.03892333837659834
$ws.Cells.Item(20, 5).Value = 0.2052821940398672
$ws.Cells.Item(20, 6).Value = 0.8272508232600302
$ws.Cells.Item(20, 11).Value = 0.3732867893952516
$ws.Cells.Item(20, 13).Value = 0.2952643798722292
$ws.Cells.Item(20, 15).Value = 2.83059448272752
$ws.Cells.Item(21, 2).Value = 0.4420203680753332
$ws.Cells.Item(21, 3).Value = 0.06765662162816
$ws.Cells.Item(21, 4).Value = 0.04123504468892492
$ws.Cells.Item(21, 5).Value = 0.2218206411589776
$ws.Cells.Item(21, 6).Value = 0.8348770491944038
$ws.Cells.Item(21, 11).Value = 0.4190978309206059
$ws.Cells.Item(21, 13).Value = 0.3267085849058304
$ws.Cells.Item(21, 15).Value = 2.836165307007263
$ws.Cells.Item(22, 2).Value = 0.4718045780902571
$ws.Cells.Item(22, 3).Value = 0.07048962052525098
$ws.Cells.Item(22, 4).Value = 0.04273584266839237
$ws.Cells.Item(22, 5).Value = 0.2327402831699317
$ws.Cells.Item(22, 6).Value = 0.8403087597600063
$ws.Cells.Item(22, 11).Value = 0.4489931986000499
$ws.Cells.Item(22, 13).Value = 0.3473298277708068
$ws.Cells.Item(22, 15).Value = 2.841287536966178
$ws.Cells.Item(23, 2).Value = 0.4559077607520123
$ws.Cells.Item(23, 3).Value = 0.0689790339592804
$ws.Cells.Item(23, 4).Value = 0.04193576594929027
$ws.Cells.Item(23, 5).Value = 0.2269020425680708
$ws.Cells.Item(23, 6).Value = 0.837368771610528
$ws.Cells.Item(23, 11).Value = 0.4330417102545994
$ws.Cells.Item(23, 13).Value = 0.3363174008265162
$ws.Cells.Item(23, 15).Value = 2.838418012816902
$ws.Cells.Item(24, 2).Value = 0.3957127114505852
$ws.Cells.Item(24, 3).Value = 0.06322657747271876
$ws.Cells.Item(24, 4).Value = 0.03888538735330371
$ws.Cells.Item(24, 5).Value = 0.2050137338956333
$ws.Cells.Item(24, 6).Value = 0.8271337244266022
$ws.Cells.Item(24, 11).Value = 0.3725373051731538
$ws.Cells.Item(24, 13).Value = 0.2947516290638603
$ws.Cells.Item(24, 15).Value = 2.830528406120408
$ws.Cells.Item(25, 2).Value = 0.3308221618610503
$ws.Cells.Item(25, 3).Value = 0.05695623117097171
$ws.Cells.Item(25, 4).Value = 0.03555288947519841
$ws.Cells.Item(25, 5).Value = 0.181870941155708
$ws.Cells.Item(25, 6).Value = 0.8180115616049903
$ws.Cells.Item(25, 11).Value = 0.3070944521543026
$ws.Cells.Item(25, 13).Value = 0.250217483831392
$ws.Cells.Item(25, 15).Value = 2.828363172449087
